$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "isRun" column: header in E1, plus one flag value per existing data row.
$ws.Range("E1").Value = "isRun"
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 0

# Mirror the author's final on-screen selection over the new column.
[void]$ws.Range("E1:E6").Select()

# Match the saved window height from the source commit.
$excel.ActiveWindow.Height = 13120
